# Update previous_count (C) and change (D) columns for agencies whose
# previous_count dropped since yesterday's comparison.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  C = 8;  D = 1 },
    @{ Row = 7;  C = 26; D = 1 },
    @{ Row = 8;  C = 22; D = 1 },
    @{ Row = 10; C = 31; D = 4 },
    @{ Row = 14; C = 25; D = 1 },
    @{ Row = 16; C = 12; D = 1 },
    @{ Row = 17; C = 28; D = 1 },
    @{ Row = 18; C = 19; D = 1 },
    @{ Row = 29; C = 15; D = 3 },
    @{ Row = 30; C = 9;  D = 1 },
    @{ Row = 35; C = 16; D = 1 },
    @{ Row = 38; C = 19; D = 1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
